# Stage 1: update companies data
# The data rows (2-34, columns A-K) are being re-ordered: each row's full
# record (name, number, dates, category, SIC codes, description, use case)
# moves as a whole to a new row position. Read the whole block once, build
# the reordered block in memory, then write it back in one shot so per-cell
# dependent formatting/formula state (none here) stays consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRange = $ws.Range("A2:K34")
$data = $srcRange.Value2

# Columns B (Company Number), C (Incorporation Date), F (Date Downloaded),
# G (Time Discovered) and I (SIC Codes) hold numeric-/date-/time-looking
# text in the source file (all stored as plain text, t="inlineStr"). Mark
# just those columns as Text before writing so Excel doesn't silently
# reinterpret e.g. "16473515" or "2025-05-26" as a number/date; leave the
# other columns alone so they keep the workbook's original (unstyled) look.
$ws.Range("B2:C34").NumberFormat = "@"
$ws.Range("F2:G34").NumberFormat = "@"
$ws.Range("I2:I34").NumberFormat = "@"

$rowCount = 33
$colCount = 11

# destRowOffset[i] = 0-based offset (within the A2:K34 block) of the source
# row that should land at destination offset i (0-based, dest row = i + 2).
$destRowOffset = @(1,0,5,3,4,2,6,10,9,7,8,13,14,11,12,18,19,16,15,17,24,25,23,21,20,22,29,28,26,27,32,31,30)

$newData = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $destRowOffset[$i] + 1
    for ($c = 1; $c -le $colCount; $c++) {
        $newData[$i, $c - 1] = $data[$srcRow, $c]
    }
}

$ws.Range("A2:K34").Value2 = $newData
